$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'68.125.45"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +0.45%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'3.765.04"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -0.78%  "

$ws.Range("E4").Value = "  +0.13%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'594.93"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.59%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'167.41"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -0.63%  "

$style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'3.760.73"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -0.69%  "

$style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.522"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -1.43%  "

$ws.Range("E10").Value = "  -2.55%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'6.48"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -0.19%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.449"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -2.30%  "

$ws.Range("E13").Value = "  -1.44%  "

$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'36.68"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -0.43%  "

$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'4.399.58"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -0.74%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'3.770.51"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -1.41%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'68.149.18"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +0.63%  "

$style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'18.20"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -3.84%  "

$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'7.07"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -3.47%  "

$ws.Range("E20").Value = "  -0.34%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'10.81"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +2.12%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'469.61"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +0.46%  "

$style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'0.703"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -3.79%  "

$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'85.04"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +1.90%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'0.0000144"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -4.82%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'2.24"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -2.63%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'12.17"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -0.22%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'10.15"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -1.49%  "

$ws.Range("E29").Value = "  +0.11%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'3.915.39"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -0.99%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'2.79"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -4.59%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'7.43"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -3.24%  "

$style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'2.25"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -0.75%  "

$style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'30.07"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -1.57%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'9.25"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +0.31%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").Value = "'0.997"
$ws.Range("D36").Style = $style

$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'3.722.17"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -0.97%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'0.102"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -3.13%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'3.48"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -6.93%  "

$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").Value = "  -0.41%  "

$style = $ws.Range("D42").Style
$ws.Range("D42").Value = "'5.81"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -2.03%  "

$ws.Range("E43").Value = "  +0.33%  "

$style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'0.308"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -2.96%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'1.94"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -1.47%  "

$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'8.59"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -1.84%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'400.05"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -1.93%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'45.52"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.65%  "

$style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'145.18"
$ws.Range("D50").Style = $style

$style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'25.57"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +0.91%  "

